{"js": "// Update the worksheet date and the 25 two-digit-by-two-digit\n// multiplication problems/answers to the new day's values.\n// Each old value is unique in the document, so an exact (case-sensitive,\n// whole-match) search-and-replace on each pair is sufficient and avoids\n// any ambiguity about which table cell / paragraph a value lives in.\n\nconst replacements = [\n  [\"2024-11-18 Monday\", \"2024-11-19 Tuesday\"],\n  [\"93\u00d785=7905\", \"92\u00d784=7728\"],\n  [\"69\u00d715=1035\", \"87\u00d754=4698\"],\n  [\"31\u00d725=775\", \"36\u00d753=1908\"],\n  [\"38\u00d750=1900\", \"32\u00d797=3104\"],\n  [\"99\u00d753=5247\", \"44\u00d712=528\"],\n  [\"29\u00d786=2494\", \"71\u00d715=1065\"],\n  [\"27\u00d786=2322\", \"44\u00d724=1056\"],\n  [\"89\u00d717=1513\", \"96\u00d722=2112\"],\n  [\"33\u00d754=1782\", \"68\u00d739=2652\"],\n  [\"18\u00d732=576\", \"16\u00d784=1344\"],\n  [\"29\u00d717=493\", \"55\u00d744=2420\"],\n  [\"54\u00d758=3132\", \"19\u00d772=1368\"],\n  [\"64\u00d796=6144\", \"90\u00d793=8370\"],\n  [\"87\u00d746=4002\", \"55\u00d771=3905\"],\n  [\"34\u00d755=1870\", \"77\u00d765=5005\"],\n  [\"42\u00d712=504\", \"84\u00d778=6552\"],\n  [\"46\u00d729=1334\", \"42\u00d779=3318\"],\n  [\"20\u00d742=840\", \"21\u00d755=1155\"],\n  [\"58\u00d783=4814\", \"30\u00d727=810\"],\n  [\"96\u00d798=9408\", \"51\u00d766=3366\"],\n  [\"69\u00d728=1932\", \"54\u00d715=810\"],\n  [\"98\u00d721=2058\", \"35\u00d734=1190\"],\n  [\"21\u00d719=399\", \"52\u00d715=780\"],\n  [\"98\u00d752=5096\", \"42\u00d733=1386\"],\n  [\"91\u00d735=3185\", \"86\u00d736=3096\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 two-digit-by-two-digit\n# multiplication problems/answers to the new day's values.\n# Each \"old\" value below is unique in the document, so a plain\n# Find/Replace (wdReplaceAll, effectively a single hit) on each pair\n# unambiguously targets the right paragraph / table cell.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-11-18 Monday\", \"2024-11-19 Tuesday\"),\n    @(\"93\u00d785=7905\", \"92\u00d784=7728\"),\n    @(\"69\u00d715=1035\", \"87\u00d754=4698\"),\n    @(\"31\u00d725=775\", \"36\u00d753=1908\"),\n    @(\"38\u00d750=1900\", \"32\u00d797=3104\"),\n    @(\"99\u00d753=5247\", \"44\u00d712=528\"),\n    @(\"29\u00d786=2494\", \"71\u00d715=1065\"),\n    @(\"27\u00d786=2322\", \"44\u00d724=1056\"),\n    @(\"89\u00d717=1513\", \"96\u00d722=2112\"),\n    @(\"33\u00d754=1782\", \"68\u00d739=2652\"),\n    @(\"18\u00d732=576\", \"16\u00d784=1344\"),\n    @(\"29\u00d717=493\", \"55\u00d744=2420\"),\n    @(\"54\u00d758=3132\", \"19\u00d772=1368\"),\n    @(\"64\u00d796=6144\", \"90\u00d793=8370\"),\n    @(\"87\u00d746=4002\", \"55\u00d771=3905\"),\n    @(\"34\u00d755=1870\", \"77\u00d765=5005\"),\n    @(\"42\u00d712=504\", \"84\u00d778=6552\"),\n    @(\"46\u00d729=1334\", \"42\u00d779=3318\"),\n    @(\"20\u00d742=840\", \"21\u00d755=1155\"),\n    @(\"58\u00d783=4814\", \"30\u00d727=810\"),\n    @(\"96\u00d798=9408\", \"51\u00d766=3366\"),\n    @(\"69\u00d728=1932\", \"54\u00d715=810\"),\n    @(\"98\u00d721=2058\", \"35\u00d734=1190\"),\n    @(\"21\u00d719=399\", \"52\u00d715=780\"),\n    @(\"98\u00d752=5096\", \"42\u00d733=1386\"),\n    @(\"91\u00d735=3185\", \"86\u00d736=3096\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $old\"\n    }\n}\n\nWrite-Output \"done\"\n"}
